$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions): update 想去人数 (attendance count) for two events
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F4").Value = 1736
$wsExpo.Range("F6").Value = 194

# Sheet "全部类型" (All types): same two events appear on different rows
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1736
$wsAll.Range("F7").Value = 194
